$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Workbook / sheet view bits
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,19).Select() | Out-Null

# Column widths (A..K)
$ws.Columns.Item(1).ColumnWidth = 53.7109375
$ws.Columns.Item(2).ColumnWidth = 20.28515625
$ws.Columns.Item(3).ColumnWidth = 47.28515625
$ws.Columns.Item(4).ColumnWidth = 25.140625
$ws.Columns.Item(5).ColumnWidth = 18.140625
$ws.Columns.Item(6).ColumnWidth = 14.7109375
$ws.Columns.Item(7).ColumnWidth = 21.5703125
$ws.Columns.Item(8).ColumnWidth = 21
$ws.Columns.Item(9).ColumnWidth = 19.85546875
$ws.Columns.Item(10).ColumnWidth = 14.42578125
$ws.Columns.Item(11).ColumnWidth = 16

# ---------------------------------------------------------------------------
# Helper cells used only to "mint" each distinct cell-style combo in the
# exact order the target workbook uses them (so the generated cellXfs table
# lines up), before being applied (and reused) on the real target cells.
# ---------------------------------------------------------------------------

# idx2: quotePrefix, otherwise default -> B28/B29/B30
$ws.Range("B28,B29,B30").Value = "'"
$ws.Range("B28,B29,B30").ClearContents()

# idx3: text number format (49) only -> L28/L29/L30/L31
$ws.Range("L28,L29,L30,L31").NumberFormat = "@"

# idx4: thin box border only, default font/fill -> B33:E33 (keeps its text)
$ws.Range("B33:E33").Borders.LineStyle = 1
$ws.Range("B33:E33").Borders.Weight = 2
$ws.Range("B33").Value = "Статический (CRTP)"
$ws.Range("C33").Value = "Статический полнофункциональный  (CRTP + meta)"
$ws.Range("D33").Value = "Внешний (смешанный)"
$ws.Range("E33").Value = "Динамический"

# idx5: "Good" style + text format + border + quotePrefix -> B34,C34
$ws.Range("B34,C34").Borders.LineStyle = 1
$ws.Range("B34,C34").Borders.Weight = 2
$ws.Range("B34,C34").Style = "Good"
$ws.Range("B34,C34").NumberFormat = "@"
$ws.Range("B34,C34").Value = "'"
$ws.Range("B34,C34").ClearContents()

# idx6: grey header fill + border, default font -> column A rows 34..45
$ws.Range("A34:A45").Borders.LineStyle = 1
$ws.Range("A34:A45").Borders.Weight = 2
$ws.Range("A34:A45").Interior.Color = 14277081
$ws.Range("A34").Value = "Память"
$ws.Range("A35").Value = "Скорость"
$ws.Range("A36").Value = "Возможность хранение в базовом типе всех наследников"
$ws.Range("A37").Value = "Добавление вирт. методов вне класса"
$ws.Range("A38").Value = "nullptr"
$ws.Range("A39").Value = "Работа со статическими типами данных"
$ws.Range("A40").Value = "Работа не только с наследниками базового"
$ws.Range("A41").Value = "Количество кода"
$ws.Range("A42").Value = "Удобство рефракторинга"
$ws.Range("A43").Value = "Простота добавления новых элементов"
$ws.Range("A44").Value = "Возможность хранить в коллекции"
$ws.Range("A45").Value = "Засорение пространства имен"

# idx7: "Good" + text format + border, no quote prefix -> E34
$ws.Range("E34").Borders.LineStyle = 1
$ws.Range("E34").Borders.Weight = 2
$ws.Range("E34").Style = "Good"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "x"
$ws.Range("E34").ClearContents()

# idx8: "Neutral" + text format + border, no quote prefix -> D35,E35
$ws.Range("D35,E35").Borders.LineStyle = 1
$ws.Range("D35,E35").Borders.Weight = 2
$ws.Range("D35,E35").Style = "Neutral"
$ws.Range("D35,E35").NumberFormat = "@"
$ws.Range("D35,E35").Value = "x"
$ws.Range("D35,E35").ClearContents()

# idx9: "Neutral" + text format + border + quote prefix -> D34
$ws.Range("D34").Borders.LineStyle = 1
$ws.Range("D34").Borders.Weight = 2
$ws.Range("D34").Style = "Neutral"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "'"
$ws.Range("D34").ClearContents()

# idx10: "Bad" + text format + border, no quote prefix -> C35
$ws.Range("C35").Borders.LineStyle = 1
$ws.Range("C35").Borders.Weight = 2
$ws.Range("C35").Style = "Bad"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "x"
$ws.Range("C35").ClearContents()

# idx11: "Bad" + border, default number format -> B39
$ws.Range("B39").Borders.LineStyle = 1
$ws.Range("B39").Borders.Weight = 2
$ws.Range("B39").Style = "Bad"

# idx12: "Good" + border, default number format -> C39
$ws.Range("C39").Borders.LineStyle = 1
$ws.Range("C39").Borders.Weight = 2
$ws.Range("C39").Style = "Good"

# idx13: "Neutral" + border, default number format -> D41
$ws.Range("D41").Borders.LineStyle = 1
$ws.Range("D41").Borders.Weight = 2
$ws.Range("D41").Style = "Neutral"

# idx14: grey fill + left/right-only border, default font -> A46
$ws.Range("A46").Borders(7).LineStyle = 1
$ws.Range("A46").Borders(7).Weight = 2
$ws.Range("A46").Borders(10).LineStyle = 1
$ws.Range("A46").Borders(10).Weight = 2
$ws.Range("A46").Interior.Color = 14277081
$ws.Range("A46").Value = "Наличие виртуальных функци"

# ---------------------------------------------------------------------------
# Now stamp every remaining matrix cell (B..E, rows 35..46) with the right
# combo (they all reuse one of the xf indices minted above).
# ---------------------------------------------------------------------------

function Set-Cell([string]$addr, [string]$kind) {
    $rng = $ws.Range($addr)
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
    switch ($kind) {
        "Good-q"    { $rng.Style = "Good";    $rng.NumberFormat = "@"; $rng.Value = "'"; $rng.ClearContents() }
        "Good"      { $rng.Style = "Good";    $rng.NumberFormat = "@"; $rng.Value = "x"; $rng.ClearContents() }
        "Bad-q"     { $rng.Style = "Bad";     $rng.NumberFormat = "@"; $rng.Value = "'"; $rng.ClearContents() }
        "Bad"       { $rng.Style = "Bad";     $rng.NumberFormat = "@"; $rng.Value = "x"; $rng.ClearContents() }
        "Neutral-q" { $rng.Style = "Neutral"; $rng.NumberFormat = "@"; $rng.Value = "'"; $rng.ClearContents() }
        "Neutral"   { $rng.Style = "Neutral"; $rng.NumberFormat = "@"; $rng.Value = "x"; $rng.ClearContents() }
        "Good-p"    { $rng.Style = "Good" }
        "Bad-p"     { $rng.Style = "Bad" }
        "Neutral-p" { $rng.Style = "Neutral" }
    }
}

# Row 36
Set-Cell "B36" "Bad"
Set-Cell "C36" "Good"
Set-Cell "D36" "Good"
Set-Cell "E36" "Good"

# Row 37
Set-Cell "B37" "Bad"
Set-Cell "C37" "Bad"
Set-Cell "D37" "Good"
Set-Cell "E37" "Bad"

# Row 38
Set-Cell "B38" "Good"
Set-Cell "C38" "Bad"
Set-Cell "D38" "Bad"
Set-Cell "E38" "Good"

# Row 39
Set-Cell "D39" "Good-p"
Set-Cell "E39" "Bad-p"

# Row 40
Set-Cell "B40" "Bad-p"
Set-Cell "C40" "Bad-p"
Set-Cell "D40" "Good-p"
Set-Cell "E40" "Bad-p"

# Row 41
Set-Cell "B41" "Good-p"
Set-Cell "C41" "Bad-p"
Set-Cell "E41" "Good-p"

# Row 42
Set-Cell "B42" "Good-p"
Set-Cell "C42" "Neutral-p"
Set-Cell "D42" "Neutral-p"
Set-Cell "E42" "Good-p"

# Row 43
Set-Cell "B43" "Good-p"
Set-Cell "C43" "Neutral-p"
Set-Cell "D43" "Good-p"
Set-Cell "E43" "Good-p"

# Row 44
Set-Cell "B44" "Bad-p"
Set-Cell "C44" "Good-p"
Set-Cell "D44" "Good-p"
Set-Cell "E44" "Good-p"

# Row 45
Set-Cell "B45" "Good-p"
Set-Cell "C45" "Neutral-p"
Set-Cell "D45" "Bad-p"
Set-Cell "E45" "Good-p"

# Row 46
Set-Cell "B46" "Bad-p"
Set-Cell "C46" "Bad-p"
Set-Cell "D46" "Good-p"
Set-Cell "E46" "Good-p"
